$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster changes from ECs(20) to FAPs, and values updated
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 1.030436666666667
$ws.Range("H2").Value = 3.09131
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.124720666666667
$ws.Range("N2").Value = 6.374162
$ws.Range("O2").Value = 0.890171494008882
$ws.Range("P2").Value = 0.8901714940088818
$ws.Range("Q2").Value = 2.189390081357778
$ws.Range("R2").Value = 19.70451073222
$ws.Range("S2").Value = 0.890171494008882
$ws.Range("T2").Value = 0.8901714940088818

# Row 3: Target cluster changes from FAPs(23) to MuSCs(24), and values updated
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 1.030436666666667
$ws.Range("H3").Value = 3.09131
$ws.Range("M3").Value = 0.262146
$ws.Range("N3").Value = 0.786438
$ws.Range("O3").Value = 0.1098285059911181
$ws.Range("P3").Value = 0.1098285059911181
$ws.Range("Q3").Value = 0.27012485042
$ws.Range("R3").Value = 2.43112365378
$ws.Range("S3").Value = 0.1098285059911181
$ws.Range("T3").Value = 0.1098285059911181

# Row 4: deleted
$ws.Range("A4:T4").Delete()
